$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (8 & 9) for the "headless" feature --------------------
# Write B8 before A8 so the shared-string table append order matches the
# target workbook (headless, start_with_code, start_manually, manual).
$ws.Range("B8").Value = "headless"
$ws.Range("A8").Value = "start_with_code"
$ws.Range("A9").Value = "start_manually"
$ws.Range("B9").Value = "manual"

# --- Formatting -------------------------------------------------------
# Column A (A8:A9) should look like the existing header cells (A1:A7),
# so copy that formatting over (reuses the existing style, no new xf).
$ws.Range("A1").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)   # xlPasteFormats

# Column B (B8:B9) starts from the same base formatting, then gets its
# own distinct style (matching the new cellXfs entry added upstream).
$ws.Range("B8:B9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B8:B9").Font.Name = "Calibri Light"
$ws.Range("B8:B9").HorizontalAlignment = -4108   # xlCenter

# --- Selection ----------------------------------------------------------
[void]$ws.Range("D6").Select()
